# 1st changes of mifos to finflux
# Insert a new blank column before column N on the "Repayment schedule" sheet.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Columns("N:N").Insert()

# Make "Repayment schedule" the active/selected sheet and set the selection.
$ws.Activate()
$ws.Range("R8").Select()
